# Center the title on the "Results (Student)" slide.
#
# The author's commit ("updates to paper and small changes to
# presentation") centers the title text box on the slide that reads
# "Results (Student)". The cached `datetimeFigureOut` field text seen
# in the slide masters/layouts of the canonical diff is a PowerPoint
# save-time artifact (today's date re-cached into the field) and isn't
# something a user/script action changes directly, so it's left alone.

$p = $ppt.ActivePresentation

foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "Results (Student)") {
                $shp.TextFrame.TextRange.ParagraphFormat.Alignment = 2  # ppAlignCenter
            }
        }
    }
}
